$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# "Nomenclature all in one table": reviewer-2 comments 2.1 and 2.2
# (rows 18-19) gain a "Response" (col C) and a "Status of changes"
# (col F) entry, matching the 6-column layout already used further
# down the sheet (rows 33-36) for the "DONE." status cells. Re-use
# the existing styles (grey-fill "done" look) by pasting formats from
# cells that already carry them, instead of inventing new styles.
# ------------------------------------------------------------------

# A18 / A19 -> style 8 (grey fill, no wrap) -- same as A33/A34.
$ws.Range("A33").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# B18 / B19 -> style 9 (grey fill, wrap) -- same as B33/B34.
$ws.Range("B33").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# D18 / D19 -> style 8 (grey fill, no wrap) -- same as D33/D34.
$ws.Range("D33").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$ws.Range("D19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# E18 / E19 -> style 8 (grey fill, no wrap) -- borrow format from a
# style-8 cell (style index is the same no matter which column it
# came from).
$ws.Range("D34").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$ws.Range("E19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# C18 -> style 8 (grey fill, no wrap) -- same as C34.
$ws.Range("C34").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# C19 -> style 9 (grey fill, wrap) -- same as C33.
$ws.Range("C33").Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# F18 / F19 -> style 8 (grey fill, no wrap) -- same as F33/F34.
$ws.Range("F33").Copy() | Out-Null
$ws.Range("F18").PasteSpecial(-4122) | Out-Null
$ws.Range("F19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# New cell content (the PasteSpecial calls above only copied format,
# so the pre-existing A/B/D/E values were left untouched).
# ------------------------------------------------------------------
$ws.Range("C18").Value = "Comma inserted, as suggested. Thanks!"
$ws.Range("C19").Value = "The section in question seems right to us. The sentence says, ""Thus, light is EM radiation … sensitivity."" The parenthetical sentence is complete as writtten. ""(See Figure 1.)"" is an imperative sentence with implied subject ""You"". No changes were made in response to this particular suggestion. However, at the reviewer's suggestion, we did another thorough proofread immediately prior to submission."
$ws.Range("F18").Value = "DONE."
$ws.Range("F19").Value = "DONE."

# Row 19 grows tall to fit the long wrapped response text.
$ws.Rows("19").RowHeight = 119

# ------------------------------------------------------------------
# Scroll/select so the newly-unified table is in view.
# ------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A18:F18").Select() | Out-Null
